$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "B4";  Value = 7.377 },
    @{ Cell = "A9";  Value = -21.74 },
    @{ Cell = "B9";  Value = 5.88 },
    @{ Cell = "D9";  Value = -7.852000000000001 },
    @{ Cell = "B11"; Value = 6.33 },
    @{ Cell = "A18"; Value = -21.995 },
    @{ Cell = "A20"; Value = -20.584 },
    @{ Cell = "B23"; Value = 7.650000000000001 },
    @{ Cell = "B24"; Value = 5.238 },
    @{ Cell = "B26"; Value = 5.595000000000001 },
    @{ Cell = "A27"; Value = -21.7 },
    @{ Cell = "D27"; Value = -7.741999999999999 },
    @{ Cell = "D29"; Value = -7.654999999999999 },
    @{ Cell = "D32"; Value = -7.258 },
    @{ Cell = "B34"; Value = 6.968999999999999 },
    @{ Cell = "A35"; Value = -21.603 },
    @{ Cell = "B35"; Value = 6.085000000000001 },
    @{ Cell = "D37"; Value = -7.636000000000001 },
    @{ Cell = "D38"; Value = -8.019 },
    @{ Cell = "D41"; Value = -7.944000000000001 },
    @{ Cell = "D45"; Value = -7.470000000000001 },
    @{ Cell = "B48"; Value = 5.433999999999999 },
    @{ Cell = "B49"; Value = 6.237 },
    @{ Cell = "D51"; Value = -8.394 },
    @{ Cell = "B52"; Value = 5.544 },
    @{ Cell = "D57"; Value = -8.145999999999999 },
    @{ Cell = "D64"; Value = -7.672 },
    @{ Cell = "B66"; Value = 5.172 },
    @{ Cell = "B67"; Value = 5.374 },
    @{ Cell = "A69"; Value = -21.376 },
    @{ Cell = "A76"; Value = -20.392 },
    @{ Cell = "A78"; Value = -20.652 },
    @{ Cell = "B78"; Value = 7.680999999999999 },
    @{ Cell = "B80"; Value = 8.301 },
    @{ Cell = "A82"; Value = -21.718 },
    @{ Cell = "D82"; Value = -8.061 },
    @{ Cell = "A83"; Value = -21.509 },
    @{ Cell = "A93"; Value = -21.392 },
    @{ Cell = "D93"; Value = -6.956 },
    @{ Cell = "B99"; Value = 5.292999999999999 },
    @{ Cell = "D102"; Value = -7.761 },
    @{ Cell = "B104"; Value = 7.057 },
    @{ Cell = "D105"; Value = -7.842000000000001 }
)

foreach ($change in $changes) {
    $ws.Range($change.Cell).Value = $change.Value
}
